# Add a "table_data" field to the KeyMeasurements sheets, placed right
# before the existing "data_reference" column, shifting that column and
# everything after it (linked_references, name, description) one to the
# right.

$wb = $excel.ActiveWorkbook

# FieldIlluminationKeyMeasurements: data_reference was in column AQ.
$ws = $wb.Worksheets.Item("FieldIlluminationKeyMeasurements")
$ws.Columns("AQ:AQ").Insert()
$ws.Range("AQ1").Value = "table_data"

# PSFBeadsKeyMeasurements: data_reference was in column BJ.
$ws = $wb.Worksheets.Item("PSFBeadsKeyMeasurements")
$ws.Columns("BJ:BJ").Insert()
$ws.Range("BJ1").Value = "table_data"
